$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3).
#    NOTE: compare with the literal string on the LEFT of -eq: some cells
#    hold real booleans, and PowerShell's -eq coerces a string RHS to bool
#    against a bool LHS (any non-empty string becomes $true), which would
#    otherwise produce false "True"/"False" matches.
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ("Ready for handoff" -eq $cell.Text) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Column widths shrink (report regenerated with narrower status columns):
#      Overview columns E (zh-cn) and F (de-de)
#      zh-cn column C (Status)
#      de-de column C (Status)
#    previously ~17.22 chars wide, now ~13.41 chars wide.
# ---------------------------------------------------------------------------
$targetWidth = 13.4101845877511 - (5.0 / 6.0)

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $targetWidth
$overview.Columns.Item(6).ColumnWidth = $targetWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $targetWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $targetWidth
